$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1159.7865
$ws.Range("I15").Value = 1159.7865
$ws.Range("K15").Value = 3479.3595
$ws.Range("M15").Value = -3310.3595
$ws.Range("H116").Value = 10594.444
$ws.Range("I116").Value = 3500
$ws.Range("J116").Value = 11481.25
$ws.Range("K116").Value = 3500
$ws.Range("L116").Value = 11481.25
$ws.Range("M116").Value = -58
$ws.Range("N116").Value = -18365.25
$ws.Range("H132").Value = 1893.175
$ws.Range("I132").Value = 1016.1111
$ws.Range("J132").Value = 3714.7693
$ws.Range("K132").Value = 3048.3333
$ws.Range("L132").Value = 11144.3079
$ws.Range("M132").Value = -518.3332999999998
$ws.Range("N132").Value = -16204.3079
$ws.Range("H135").Value = 367.44446
$ws.Range("I135").Value = 367.44446
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3307.00014
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -772.0001400000001
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1895.1111
$ws.Range("I137").Value = 1321.2
$ws.Range("J137").Value = 2115.8462
$ws.Range("K137").Value = 3963.6
$ws.Range("L137").Value = 6347.5386
$ws.Range("M137").Value = -1413.6
$ws.Range("N137").Value = -11447.5386
$ws.Range("H138").Value = 1591.5161
$ws.Range("I138").Value = 785.6326
$ws.Range("J138").Value = 4629.077
$ws.Range("K138").Value = 2356.8978
$ws.Range("L138").Value = 13887.231
$ws.Range("M138").Value = 2783.1022
$ws.Range("N138").Value = -24167.231
$ws.Range("H141").Value = 6029.9287
$ws.Range("I141").Value = 4353.625
$ws.Range("J141").Value = 8265
$ws.Range("K141").Value = 13060.875
$ws.Range("L141").Value = 24795
$ws.Range("M141").Value = -7880.875
$ws.Range("N141").Value = -35155
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1855.091
$ws.Range("I61").Value = 1565.3334
$ws.Range("J61").Value = 2202.8
$ws.Range("K61").Value = 1565.3334
$ws.Range("L61").Value = 2202.8
$ws.Range("M61").Value = -1353.3334
$ws.Range("N61").Value = -2626.8
$ws.Range("H97").Value = 788.8378
$ws.Range("I97").Value = 635.0645
$ws.Range("J97").Value = 1583.3334
$ws.Range("K97").Value = 635.0645
$ws.Range("L97").Value = 1583.3334
$ws.Range("M97").Value = -139.0645
$ws.Range("N97").Value = -2575.3334
$ws.Range("H132").Value = 1854.3846
$ws.Range("I132").Value = 1285.4166
$ws.Range("J132").Value = 3134.5625
$ws.Range("K132").Value = 3856.2498
$ws.Range("L132").Value = 9403.6875
$ws.Range("M132").Value = -1326.2498
$ws.Range("N132").Value = -14463.6875
$ws.Range("H136").Value = 1855.091
$ws.Range("I136").Value = 1565.3334
$ws.Range("J136").Value = 2202.8
$ws.Range("K136").Value = 4696.0002
$ws.Range("L136").Value = 6608.400000000001
$ws.Range("M136").Value = -2146.0002
$ws.Range("N136").Value = -11708.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3143.6
$ws.Range("I31").Value = 1964.421
$ws.Range("J31").Value = 4543.875
$ws.Range("K31").Value = 1964.421
$ws.Range("L31").Value = 4543.875
$ws.Range("M31").Value = -1669.421
$ws.Range("N31").Value = -5133.875
$ws.Range("H34").Value = 3143.6
$ws.Range("I34").Value = 1964.421
$ws.Range("J34").Value = 4543.875
$ws.Range("K34").Value = 1964.421
$ws.Range("L34").Value = 4543.875
$ws.Range("M34").Value = -1762.421
$ws.Range("N34").Value = -4947.875
$ws.Range("H58").Value = 2580
$ws.Range("I58").Value = 1539.5555
$ws.Range("J58").Value = 3300.3076
$ws.Range("K58").Value = 1539.5555
$ws.Range("L58").Value = 3300.3076
$ws.Range("M58").Value = -1336.5555
$ws.Range("N58").Value = -3706.3076
$ws.Range("H132").Value = 1578.0555
$ws.Range("I132").Value = 1078.2609
$ws.Range("J132").Value = 2462.3076
$ws.Range("K132").Value = 3234.7827
$ws.Range("L132").Value = 7386.9228
$ws.Range("M132").Value = -704.7826999999997
$ws.Range("N132").Value = -12446.9228
$ws.Range("H136").Value = 2580
$ws.Range("I136").Value = 1539.5555
$ws.Range("J136").Value = 3300.3076
$ws.Range("K136").Value = 4618.666499999999
$ws.Range("L136").Value = 9900.9228
$ws.Range("M136").Value = -2068.666499999999
$ws.Range("N136").Value = -15000.9228
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1681.25
$ws.Range("J34").Value = 2183.3333
$ws.Range("L34").Value = 6549.999899999999
$ws.Range("N34").Value = -6717.999899999999
$ws.Range("H39").Value = 4526.9443
$ws.Range("H55").Value = 39230.76
$ws.Range("I55").Value = 65403.832
$ws.Range("J55").Value = 4333.3335
$ws.Range("K55").Value = 196211.496
$ws.Range("L55").Value = 13000.0005
$ws.Range("M55").Value = -196034.496
$ws.Range("N55").Value = -13354.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 9800
$ws.Range("I22").Value = 9800
$ws.Range("K22").Value = 9800
$ws.Range("M22").Value = -9271
$ws.Range("H97").Value = 1393.6666
$ws.Range("I97").Value = 1134.25
$ws.Range("J97").Value = 1912.5
$ws.Range("K97").Value = 1134.25
$ws.Range("L97").Value = 1912.5
$ws.Range("M97").Value = -638.25
$ws.Range("N97").Value = -2904.5
$ws.Range("H132").Value = 5935.9653
$ws.Range("I132").Value = 7335.9443
$ws.Range("J132").Value = 3645.0908
$ws.Range("K132").Value = 22007.8329
$ws.Range("L132").Value = 10935.2724
$ws.Range("M132").Value = -19477.8329
$ws.Range("N132").Value = -15995.2724
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 604.875
$ws.Range("J22").Value = 750
$ws.Range("L22").Value = 750
$ws.Range("N22").Value = -1340
$ws.Range("H27").Value = 604.875
$ws.Range("J27").Value = 750
$ws.Range("L27").Value = 750
$ws.Range("N27").Value = -964
$ws.Range("H136").Value = 2959
$ws.Range("I136").Value = 786.8570999999999
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 2360.5713
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = 189.4287000000004
$ws.Range("N136").Value = -23100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 17375
$ws.Range("J30").Value = 17375
$ws.Range("L30").Value = 17375
$ws.Range("N30").Value = -17589
$ws.Range("H69").Value = 15973.143
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 16968.666
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 16968.666
$ws.Range("M69").Value = -9251
$ws.Range("N69").Value = -18466.666
$ws.Range("H72").Value = 15973.143
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 16968.666
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 50905.99800000001
$ws.Range("M72").Value = -26256
$ws.Range("N72").Value = -58393.99800000001
$ws.Range("H132").Value = 1559.2051
$ws.Range("I132").Value = 981.9259
$ws.Range("J132").Value = 2858.0833
$ws.Range("K132").Value = 2945.7777
$ws.Range("L132").Value = 8574.249899999999
$ws.Range("M132").Value = -415.7776999999996
$ws.Range("N132").Value = -13634.2499
$ws.Range("H136").Value = 2013.4762
$ws.Range("I136").Value = 1630.6666
$ws.Range("J136").Value = 2523.889
$ws.Range("L136").Value = 7571.667
$ws.Range("M136").Value = -2341.9998
$ws.Range("N136").Value = -12671.667
